$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '65.949.61'
$ws.Cells.Item(2, 5).Value = '  +0.14%  '

$ws.Cells.Item(3, 4).Value = '3.367.76'
$ws.Cells.Item(3, 5).Value = '  -3.36%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '581.49'
$ws.Cells.Item(5, 5).Value = '  -0.25%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '179.45'
$ws.Cells.Item(6, 5).Value = '  +3.04%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.625'
$ws.Cells.Item(7, 5).Value = '  +4.62%  '

$ws.Cells.Item(8, 5).Value = '  +0.05%  '

$ws.Cells.Item(9, 4).Value = '3.361.27'
$ws.Cells.Item(9, 5).Value = '  -3.45%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.130'
$ws.Cells.Item(10, 5).Value = '  -0.41%  '

$ws.Cells.Item(11, 5).Value = '  +1.22%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.414'
$ws.Cells.Item(12, 5).Value = '  +0.82%  '

$ws.Cells.Item(13, 4).Value = '3.962.73'
$ws.Cells.Item(13, 5).Value = '  -3.02%  '

$ws.Cells.Item(14, 5).Value = '  +0.85%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '28.90'
$ws.Cells.Item(15, 5).Value = '  -3.97%  '

$ws.Cells.Item(16, 4).Value = '65.964.25'
$ws.Cells.Item(16, 5).Value = '  -0.05%  '

$ws.Cells.Item(17, 5).Value = '  -0.46%  '

$ws.Cells.Item(18, 4).Value = '3.359.21'
$ws.Cells.Item(18, 5).Value = '  -3.60%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '5.79'
$ws.Cells.Item(19, 5).Value = '  -2.64%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '13.61'
$ws.Cells.Item(20, 5).Value = '  -2.15%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '365.52'
$ws.Cells.Item(21, 5).Value = '  -0.14%  '

$ws.Cells.Item(22, 5).Value = '  -3.02%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '72.23'
$ws.Cells.Item(23, 5).Value = '  -0.36%  '

$ws.Cells.Item(24, 5).Value = '  -0.25%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.528'
$ws.Cells.Item(25, 5).Value = '  -1.36%  '

$ws.Cells.Item(26, 5).Value = '  -0.89%  '

$ws.Cells.Item(27, 5).Value = '  +0.81%  '

$ws.Cells.Item(28, 5).Value = '  +0.43%  '

$ws.Cells.Item(29, 5).Value = '  +0.26%  '

$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.98'
$ws.Cells.Item(30, 5).Value = '  -0.28%  '

$ws.Cells.Item(31, 2).Value = 'NEARProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '5.73'
$ws.Cells.Item(31, 5).Value = '  -0.74%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '23.05'
$ws.Cells.Item(32, 5).Value = '  -4.42%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.999'
$ws.Cells.Item(33, 5).Value = '  -0.03%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '6.96'
$ws.Cells.Item(34, 5).Value = '  -2.63%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.23'
$ws.Cells.Item(35, 5).Value = '  -4.18%  '

$ws.Cells.Item(36, 5).Value = '  -1.51%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '160.88'
$ws.Cells.Item(37, 5).Value = '  +0.51%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.853'
$ws.Cells.Item(38, 5).Value = '  -3.91%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '27.27'
$ws.Cells.Item(39, 5).Value = '  -7.74%  '

$ws.Cells.Item(40, 5).Value = '  +0.42%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.61'
$ws.Cells.Item(41, 5).Value = '  +1.44%  '

$ws.Cells.Item(42, 2).Value = 'Maker'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42, 4).Value = '2.687.00'
$ws.Cells.Item(42, 5).Value = '  -4.90%  '

$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '4.33'
$ws.Cells.Item(43, 5).Value = '  -2.57%  '

$ws.Cells.Item(44, 2).Value = 'RenderToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '6.27'
$ws.Cells.Item(44, 5).Value = '  -2.32%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0674'
$ws.Cells.Item(45, 5).Value = '  -1.43%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '338.15'
$ws.Cells.Item(46, 5).Value = '  +9.96%  '

$ws.Cells.Item(47, 5).Value = '  +0.01%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '24.49'
$ws.Cells.Item(48, 5).Value = '  +0.98%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0282'
$ws.Cells.Item(49, 5).Value = '  -1.91%  '

$ws.Cells.Item(50, 5).Value = '  +3.37%  '

$ws.Cells.Item(51, 5).Value = '  +0.51%  '
